$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '44.093.49'
$ws.Cells.Item(2, 5).Value = '  +4.50%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.258.43'
$ws.Cells.Item(3, 5).Value = '  +1.97%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.24%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '229.75'
$ws.Cells.Item(5, 5).Value = '  -0.62%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.631'
$ws.Cells.Item(6, 5).Value = '  +2.31%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '63.32'
$ws.Cells.Item(7, 5).Value = '  +4.17%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.14%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.445'
$ws.Cells.Item(9, 5).Value = '  +10.65%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.103'
$ws.Cells.Item(10, 5).Value = '  +14.83%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '56.92'
$ws.Cells.Item(11, 5).Value = '  -0.70%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '26.24'
$ws.Cells.Item(12, 5).Value = '  +18.60%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.18%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.605.56'
$ws.Cells.Item(14, 5).Value = '  +2.51%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '15.63'
$ws.Cells.Item(15, 5).Value = '  +1.25%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '6.14'
$ws.Cells.Item(16, 5).Value = '  +10.12%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.832'
$ws.Cells.Item(17, 5).Value = '  +4.53%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.264.71'
$ws.Cells.Item(18, 5).Value = '  +1.50%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '44.020.27'
$ws.Cells.Item(19, 5).Value = '  +4.56%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0000101'
$ws.Cells.Item(20, 5).Value = '  +8.11%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '73.20'
$ws.Cells.Item(21, 5).Value = '  +1.76%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.01'
$ws.Cells.Item(22, 5).Value = '  -2.65%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '250.92'
$ws.Cells.Item(23, 5).Value = '  +3.20%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.06%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.41'
$ws.Cells.Item(25, 5).Value = '  -0.48%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -1.78%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '3.28'
$ws.Cells.Item(27, 5).Value = '  +23.00%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.96'
$ws.Cells.Item(28, 5).Value = '  +4.12%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '171.86'
$ws.Cells.Item(29, 5).Value = '  +1.72%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '20.70'
$ws.Cells.Item(30, 5).Value = '  +2.05%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.136'
$ws.Cells.Item(31, 5).Value = '  -2.59%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -5.45%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +2.58%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0680'
$ws.Cells.Item(34, 5).Value = '  +5.13%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.71'
$ws.Cells.Item(35, 5).Value = '  +2.42%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.85'
$ws.Cells.Item(36, 5).Value = '  -2.39%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.79'
$ws.Cells.Item(37, 5).Value = '  +6.77%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '6.62'
$ws.Cells.Item(38, 5).Value = '  +5.23%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.30'
$ws.Cells.Item(39, 5).Value = '  -1.31%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0257'
$ws.Cells.Item(40, 5).Value = '  +3.82%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.20%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '8.29'
$ws.Cells.Item(42, 5).Value = '  -2.89%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '17.40'
$ws.Cells.Item(43, 5).Value = '  +8.36%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0962'
$ws.Cells.Item(44, 5).Value = '  +0.72%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '97.28'
$ws.Cells.Item(45, 5).Value = '  +0.42%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -1.06%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'FTXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '4.33'
$ws.Cells.Item(47, 5).Value = '  -0.36%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'TerraClassic'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.000208'
$ws.Cells.Item(48, 5).Value = '  -8.51%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '1.438.50'
$ws.Cells.Item(49, 5).Value = '  -1.11%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.28'
$ws.Cells.Item(50, 5).Value = '  +3.80%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'HuobiToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.74'
$ws.Cells.Item(51, 5).Value = '  +1.40%  '
